# "modify year option to all forms"
# Insert a new column before column B ("比賽年份 Year of Competition"),
# shifting the existing header columns (old B..R) one place right (new C..S).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at B; everything from old column B onward
# (Group Name .. Team Status) shifts right by one column.
$ws.Columns("B").Insert()

# New column B gets the same width as column A (both end up at width 12).
$ws.Columns("B").ColumnWidth = $ws.Columns("A").ColumnWidth

# Populate the new header cell with the new shared string.
$ws.Range("B1").Value = "比賽年份 Year of Competition"

# Match the resulting active-cell selection recorded in the saved file.
$ws.Range("A2").Select()
